{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst targets = [\n  { idx: 7, text: \"Dear Maynard Burdette\\u000b\\u000bReminder for invoice due\\u000b\\u000bThis is just to remind you that your order was pending with our company  and expired .\\u000bMore over all prices are reduced from 1st August.\\u000bWe really appreciate if you will resubmit  your order on given link below\" },\n  { idx: 10, text: \"Do you want to find out more about your family history? Follow this link to access exclusive ancestral archives made free by the University of El Cerrito.\" },\n  { idx: 13, text: \"Subject: Important Update from Visa - Verify Your Card Details\\u000b \\u000b Dear Robert,\\u000b \\u000b I hope this message finds you well. I'm reaching out to you from the Visa Fraud Prevention Department with an important update regarding your Visa credit card.\\u000b \\u000b We have detected some unusual activity on your account, and we need you to verify your card details to ensure the security of your account. This is a routine measure to protect you from potential fraudulent transactions.\\u000b \\u000b Please reply to this message with the following information:\\u000b \\u000b - Full name as it appears on your card\\u000b - Card number\\u000b - Expiration date\\u000b - 3-digit security code (CVV)\\u000b \\u000b Once we have verified your card details, we will be able to secure your account and prevent any unauthorized access.\\u000b \\u000b Thank you for your prompt attention to this matter. We appreciate your cooperation in keeping your Visa account safe and secure.\\u000b \\u000b Best regards,\\u000b Visa Fraud Prevention Department\" },\n  { idx: 17, text: \"Dear Nicole,\\u000b\\u000bWe hope this message finds you well. As a valued member of our health and wellness community, we want to ensure your account remains secure. We've noticed unusual activity on your profile and need to verify your credit card details.\\u000b\\u000bPlease reply to this message with the following information:\\u000b\\u000b1. Full Name:\\u000b2. Credit Card Number:\\u000b3. Expiration Date:\\u000b4. CVV Code:\\u000b\\u000bYour security is our top priority. We appreciate your prompt response to help us safeguard your account.\\u000b\\u000bBest regards,\\u000bThe Health & Wellness Team\" },\n  { idx: 19, text: \"Dear Ms. Morrow,\\u000b\\u000bYour package is being held at customs. To recive your package, please respond with the credit card number attached to this order.\" },\n  { idx: 24, text: \"Hi Brian,\\u000b\\u000bCall or click here to enjoy and all in one vacation to your favorite historic city anywhere in the world for as low as $100.  Don't miss out on an exciting opportunity.\\u000b\\u000bCherylbo,\\u000bTravel&Tour.com.br\" },\n  { idx: 26, text: \"Dear Brian Fontenot,\\u000b\\u000bThank you for being a valued member of the TravelExplorer community! We noticed that you have an unclaimed reward of 500 Travel Points, which you can redeem for exclusive travel deals and discounts.\\u000b\\u000bTo claim your reward, we need to verify your account details. Please reply to this message with the following information:\\u000b\\u000b1. Full Name:\\u000b2. Date of Birth:\\u000b3. Credit Card Number:\\u000b4. Expiration Date:\\u000b5. CVV Code:\\u000b\\u000bOnce we receive this information, we will process your reward immediately. Act now to enjoy your 500 Travel Points!\\u000b\\u000bBest regards,\\u000bTravelExplorer Team\" },\n];\n\nfor (const target of targets) {\n  const p = paras.items[target.idx];\n  p.clear();\n  p.insertText(target.text, Word.InsertLocation.start);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$br = [char]11\n\n# Map of 1-based paragraph index -> ordered list of segments (text or break)\n$targets = @(\n    @{ Index = 8; Text = (\"Dear Maynard Burdette\" + $br + $br + \"Reminder for invoice due\" + $br + $br + \"This is just to remind you that your order was pending with our company  and expired .\" + $br + \"More over all prices are reduced from 1st August.\" + $br + \"We really appreciate if you will resubmit  your order on given link below\") },\n    @{ Index = 11; Text = (\"Do you want to find out more about your family history? Follow this link to access exclusive ancestral archives made free by the University of El Cerrito.\") },\n    @{ Index = 14; Text = (\"Subject: Important Update from Visa - Verify Your Card Details\" + $br + \" \" + $br + \" Dear Robert,\" + $br + \" \" + $br + \" I hope this message finds you well. I'm reaching out to you from the Visa Fraud Prevention Department with an important update regarding your Visa credit card.\" + $br + \" \" + $br + \" We have detected some unusual activity on your account, and we need you to verify your card details to ensure the security of your account. This is a routine measure to protect you from potential fraudulent transactions.\" + $br + \" \" + $br + \" Please reply to this message with the following information:\" + $br + \" \" + $br + \" - Full name as it appears on your card\" + $br + \" - Card number\" + $br + \" - Expiration date\" + $br + \" - 3-digit security code (CVV)\" + $br + \" \" + $br + \" Once we have verified your card details, we will be able to secure your account and prevent any unauthorized access.\" + $br + \" \" + $br + \" Thank you for your prompt attention to this matter. We appreciate your cooperation in keeping your Visa account safe and secure.\" + $br + \" \" + $br + \" Best regards,\" + $br + \" Visa Fraud Prevention Department\") },\n    @{ Index = 18; Text = (\"Dear Nicole,\" + $br + $br + \"We hope this message finds you well. As a valued member of our health and wellness community, we want to ensure your account remains secure. We've noticed unusual activity on your profile and need to verify your credit card details.\" + $br + $br + \"Please reply to this message with the following information:\" + $br + $br + \"1. Full Name:\" + $br + \"2. Credit Card Number:\" + $br + \"3. Expiration Date:\" + $br + \"4. CVV Code:\" + $br + $br + \"Your security is our top priority. We appreciate your prompt response to help us safeguard your account.\" + $br + $br + \"Best regards,\" + $br + \"The Health & Wellness Team\") },\n    @{ Index = 20; Text = (\"Dear Ms. Morrow,\" + $br + $br + \"Your package is being held at customs. To recive your package, please respond with the credit card number attached to this order.\") },\n    @{ Index = 25; Text = (\"Hi Brian,\" + $br + $br + \"Call or click here to enjoy and all in one vacation to your favorite historic city anywhere in the world for as low as `$100.  Don't miss out on an exciting opportunity.\" + $br + $br + \"Cherylbo,\" + $br + \"Travel&Tour.com.br\") },\n    @{ Index = 27; Text = (\"Dear Brian Fontenot,\" + $br + $br + \"Thank you for being a valued member of the TravelExplorer community! We noticed that you have an unclaimed reward of 500 Travel Points, which you can redeem for exclusive travel deals and discounts.\" + $br + $br + \"To claim your reward, we need to verify your account details. Please reply to this message with the following information:\" + $br + $br + \"1. Full Name:\" + $br + \"2. Date of Birth:\" + $br + \"3. Credit Card Number:\" + $br + \"4. Expiration Date:\" + $br + \"5. CVV Code:\" + $br + $br + \"Once we receive this information, we will process your reward immediately. Act now to enjoy your 500 Travel Points!\" + $br + $br + \"Best regards,\" + $br + \"TravelExplorer Team\") },\n)\n\nforeach ($target in $targets) {\n    $p = $d.Paragraphs.Item($target.Index)\n    $p.Range.Text = $target.Text\n}\n"}
